# Add execution estimate date for tasks - adds "TM Review" and "Sheet1" sheets,
# plus updates the current selection/view on the existing sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Evaulation 2024 - just move the selection
# ---------------------------------------------------------------------------
$wsEval = $wb.Worksheets.Item(1)
$wsEval.Activate()
$wsEval.Range("N43").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. goals 2025 - move the selection
# ---------------------------------------------------------------------------
$wsGoals = $wb.Worksheets.Item(2)
$wsGoals.Activate()
$wsGoals.Range("N10:N13").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Objectives - move the selection
# ---------------------------------------------------------------------------
$wsObjectives = $wb.Worksheets.Item(3)
$wsObjectives.Activate()
$wsObjectives.Range("A21:G21").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. PI planning - move the selection (tabSelected will move off this sheet
#    automatically once a later sheet is activated below)
# ---------------------------------------------------------------------------
$wsPI = $wb.Worksheets.Item(4)
$wsPI.Activate()
$wsPI.Range("C53:O54").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. New sheet "TM Review"
# ---------------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$wsReview = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($lastIndex))
$wsReview.Name = "TM Review"

$wsReview.Range("D5").Value = "Task management Review"

$wsReview.Range("D7").Value = "View"
$wsReview.Range("E7").Value = "Everything Is ok"

$wsReview.Range("D9").Value = "Add activity"
$wsReview.Range("E11").Value = "put filters in frame"
$wsReview.Range("E12").Value = "try colapse all stories as default view"
$wsReview.Range("E10").Value = "Select sprint shows only sprint from PI"
$wsReview.Range("E9").Value = "Add PI dropdown with all option"

$wsReview.Range("D14").Value = "Edit"
$wsReview.Range("E14").Value = "select epic add all option as default"
$wsReview.Range("E16").Value = "Add ability to create story and assign to backlog"
$wsReview.Range("E15").Value = "add filter with current sprint and current pi"
$wsReview.Range("E17").Value = "add ability to change epic for story"
$wsReview.Range("E18").Value = "add to story allighment to Objectives"

$wsReview.Range("D21").Value = "Task"
$wsReview.Range("E21").Value = "Add execution estimate date to tasks"
$wsReview.Range("E22").Value = "Create simple  notiffication email (using python anywhere)"

$wsReview.Columns.Item(4).ColumnWidth = 21.166666666666668
$wsReview.Columns.Item(5).ColumnWidth = 33.666666666666664

$wsReview.Range("J20").Select() | Out-Null

# ---------------------------------------------------------------------------
# 6. New sheet "Sheet1"
# ---------------------------------------------------------------------------
$lastIndex2 = $wb.Worksheets.Count
$wsPlan = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($lastIndex2))
$wsPlan.Name = "Sheet1"

$wsPlan.Range("F12").Value = "plan glowny"
$wsPlan.Range("G13").Value = "wyjscie na pilsko"
$wsPlan.Range("G14").Value = "prezet"
$wsPlan.Range("G15").Value = "kolacja w piątek"
$wsPlan.Range("F19").Value = "w razie nie pogody"
$wsPlan.Range("G21").Value = "lodowisko"
$wsPlan.Range("G20").Value = "śniadanie"
$wsPlan.Range("G22").Value = "kino"
$wsPlan.Range("G23").Value = "obiad"

$wsPlan.Columns.Item(6).ColumnWidth = 14.833333333333332
$wsPlan.Columns.Item(17).ColumnWidth = 14.333333333333332

$wsPlan.Range("N31").Select() | Out-Null
$wsPlan.Activate()
